# Commit before changing to progressive PrestigeUpgrades method
#
# - prod_multi!D3  : 5.0000000000000002E+25  -> 5000000000  (downstream D4:D40 formulas recalc)
# - demand_multi!C4: 5.0000000000000002E+25  -> 9000000000  (downstream C5:C41 formulas recalc)
# - selection / active-sheet bookkeeping updated to match where the author ended up working

$wb = $excel.ActiveWorkbook

# --- Update the two base values that the rest of each table cascades from ---
$wsProdMulti = $wb.Worksheets.Item("prod_multi")
$wsProdMulti.Range("D3").Value = 5000000000

$wsDemandMulti = $wb.Worksheets.Item("demand_multi")
$wsDemandMulti.Range("C4").Value = 9000000000

# --- Leave a selection on demand_multi at C5, without it staying the active tab ---
[void]$wsDemandMulti.Activate()
$wsDemandMulti.Range("C5").Select() | Out-Null

# --- End with prod_multi active and D4 selected (the final, saved UI state) ---
[void]$wsProdMulti.Activate()
$wsProdMulti.Range("D4").Select() | Out-Null
